# Auto-generated edit script: update cryptos list with refreshed prices/volumes
# and reorder a few rows whose underlying ranking shifted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.415.95"
$ws.Range("E2").Value = "'  +1.58%  "
$ws.Range("D3").Value = "'2.158.38"
$ws.Range("E3").Value = "'  +3.04%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'228.15"
$ws.Range("E5").Value = "'  -0.44%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "'  +0.49%  "
$ws.Range("D7").Value = "'64.20"
$ws.Range("E7").Value = "'  +4.59%  "
$ws.Range("E8").Value = "'  +0.03%  "
$ws.Range("E9").Value = "'  +2.81%  "
$ws.Range("E10").Value = "'  +1.96%  "
$ws.Range("E11").Value = "'  -0.14%  "
$ws.Range("D12").Value = "'15.98"
$ws.Range("E12").Value = "'  +4.09%  "
$ws.Range("D13").Value = "'2.479.64"
$ws.Range("E13").Value = "'  +3.12%  "
$ws.Range("D14").Value = "'22.33"
$ws.Range("E14").Value = "'  +1.08%  "
$ws.Range("D15").Value = "'0.813"
$ws.Range("E15").Value = "'  +0.85%  "
$ws.Range("D16").Value = "'5.54"
$ws.Range("E16").Value = "'  +0.89%  "
$ws.Range("D17").Value = "'2.157.26"
$ws.Range("E17").Value = "'  +2.94%  "
$ws.Range("D18").Value = "'39.383.77"
$ws.Range("E18").Value = "'  +1.66%  "
$ws.Range("D19").Value = "'71.84"
$ws.Range("E19").Value = "'  +0.04%  "
$ws.Range("E20").Value = "'  +0.51%  "
$ws.Range("E21").Value = "'  +1.58%  "
$ws.Range("D22").Value = "'231.41"
$ws.Range("E22").Value = "'  +1.50%  "
$ws.Range("E23").Value = "'  +0.03%  "
$ws.Range("D24").Value = "'2.54"
$ws.Range("E24").Value = "'  +7.28%  "
$ws.Range("E25").Value = "'  +0.51%  "
$ws.Range("D26").Value = "'172.21"
$ws.Range("E26").Value = "'  +0.39%  "
$ws.Range("E27").Value = "'  -0.19%  "
$ws.Range("E28").Value = "'  +0.94%  "
$ws.Range("D29").Value = "'19.93"
$ws.Range("E29").Value = "'  +3.11%  "
$ws.Range("E30").Value = "'  -1.29%  "
$ws.Range("E31").Value = "'  +8.80%  "
$ws.Range("E32").Value = "'  +1.22%  "
$ws.Range("E33").Value = "'  +2.11%  "
$ws.Range("E34").Value = "'  +0.25%  "
$ws.Range("E36").Value = "'  +1.22%  "
$ws.Range("E37").Value = "'  +0.23%  "
$ws.Range("D38").Value = "'3.60"
$ws.Range("E38").Value = "'  -0.03%  "
$ws.Range("E39").Value = "'  +0.13%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'103.78"
$ws.Range("E40").Value = "'  +2.68%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0231"
$ws.Range("E41").Value = "'  +1.02%  "
$ws.Range("D42").Value = "'17.83"
$ws.Range("E42").Value = "'  -1.16%  "
$ws.Range("D43").Value = "'1.540.03"
$ws.Range("E43").Value = "'  +0.24%  "
$ws.Range("E44").Value = "'  +4.37%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'7.92"
$ws.Range("E45").Value = "'  +3.23%  "
$ws.Range("E46").Value = "'  +0.56%  "
$ws.Range("E47").Value = "'  +1.43%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.10"
$ws.Range("E48").Value = "'  +5.87%  "
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").Value = "'4.24"
$ws.Range("E49").Value = "'  +2.83%  "
$ws.Range("D50").Value = "'2.363.57"
$ws.Range("E50").Value = "'  +3.17%  "
$ws.Range("E51").Value = "'  -0.22%  "
